# Update max_flow graph workbook:
#  - Vertices sheet (sheet1): rework vertex annotations (sources/sinks
#    gain explicit numeric excess values, edge-capacity strings shrink)
#  - Edges sheet (sheet2): drop the reverse/residual edges (rows 11-19),
#    update edge count in A1
#  - Selections updated to match where the author last clicked

$wb = $excel.ActiveWorkbook

$wsVertices = $wb.Worksheets.Item(1)
$wsEdges = $wb.Worksheets.Item(2)

# --- Vertices sheet -------------------------------------------------
# B2:B4 keep referencing the same rows, just the encoded string shrinks
# (the trailing ",<weight>" used for the old MST edge list is gone, or a
# source gets an explicit ",-1" excess tag).
$wsVertices.Range("B2").Value = "1,8,4,-1"
$wsVertices.Range("B3").Value = "3,3"
$wsVertices.Range("B4").Value = "8,10"

# B5/B6 switch from the old "parent,dist" strings to plain numeric
# capacities; B7/B8 no longer carry a value at all; B9 becomes a plain
# negative excess number.
$wsVertices.Range("B5").Value = 7
$wsVertices.Range("B6").Value = 3
$wsVertices.Range("B7").ClearContents()
$wsVertices.Range("B8").ClearContents()
$wsVertices.Range("B9").Value = -2

# --- Edges sheet ------------------------------------------------------
$wsEdges.Range("A1").Value = 9
$wsEdges.Range("A11:B19").ClearContents()

# --- Selections (Edges first so Vertices ends up the active tab) ------
$wsEdges.Range("F7").Select()
$wsVertices.Range("B7").Select()
